$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    # Force the cell to be stored as text (shared string) even when the
    # text looks like a number, without leaving a lingering custom style
    # on the cell (reset back to the Normal style afterwards).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Restricciones_del_follower sheet (index 3) ---
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = "5.550000000000001 - 3x + 1.1102230246251565e-16y"
Set-TextValue $ws.Range("B2") "-2.5500000000000003"
Set-TextValue $ws.Range("D2") "0.13"
Set-TextValue $ws.Range("E2") "5.8"
Set-TextValue $ws.Range("F2") "0"

$ws.Range("A3").Value = "0.2657500000000006 + x - 0.455y"
Set-TextValue $ws.Range("B3") "-4.265750000000001"
Set-TextValue $ws.Range("D3") "0.21"
Set-TextValue $ws.Range("E3") "4.1"
Set-TextValue $ws.Range("F3") "4.8"

$ws.Range("A4").Value = "-12.149999999999999 + x + 1.1102230246251565e-16y"
Set-TextValue $ws.Range("B4") "-5.149999999999999"
Set-TextValue $ws.Range("D4") "0.91"
Set-TextValue $ws.Range("E4") "5.6000000000000005"
Set-TextValue $ws.Range("F4") "0"

# --- Punto_modificado sheet (index 4) ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "1.85"
Set-TextValue $ws.Range("B2") "4.65"

# --- Vector_bf sheet (index 5) ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "-4.42945"

# --- Vector_BF sheet (index 6) ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "14.134999999999998"
Set-TextValue $ws.Range("A3") "-39.514500000000005"

# --- Vector_Alpha sheet (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("A2").Value = 0.09
